$wb = $excel.ActiveWorkbook
$s1 = $wb.Worksheets.Item(1)
$s1.Name = "Project Information"
$s2 = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $s1)
$s2.Name = "Authentication"
$s3 = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $s2)
$s3.Name = "Authorization"
$s2.Activate()

$s2.Range("C3").NumberFormat = "m/d/yyyy"
$s2.Range("C3").Value = 46056
